# Add a new data row (row 7) to the report sheet for the kNN method
# run with 250000 samples and k=12, matching the formatting of the
# existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 6) down to the
# new row so the new row keeps the same borders/alignment/style.
$ws.Range("A6:G6").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false | Out-Null

# Fill in the new row's values.
$ws.Range("A7").Value = "kNN"
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 250000
$ws.Range("D7").Value = "CV, 5"
$ws.Range("E7").Value = 0.81
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 920.64

# Move the active selection to the cell below the new row, mirroring the
# original workbook's "next empty cell" selection behavior.
$ws.Range("G8").Select() | Out-Null
